$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 82.83048866666665
$ws.Range("H2").Value = 248.491466
$ws.Range("I2").Value = 0.3167437020391103
$ws.Range("J2").Value = 0.3167437020391103
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.77057966666666
$ws.Range("N2").Value = 71.31173899999999
$ws.Range("O2").Value = 0.3626243450559418
$ws.Range("P2").Value = 0.3626243450559418
$ws.Range("Q2").Value = 1968.92872967993
$ws.Range("R2").Value = 17720.35856711937
$ws.Range("S2").Value = 0.1148589775025267
$ws.Range("T2").Value = 0.1148589775025268

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 82.83048866666665
$ws.Range("H3").Value = 248.491466
$ws.Range("I3").Value = 0.3167437020391103
$ws.Range("J3").Value = 0.3167437020391103
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 29.46642766666666
$ws.Range("N3").Value = 88.399283
$ws.Range("O3").Value = 0.4495155012457325
$ws.Range("P3").Value = 0.4495155012457325
$ws.Range("Q3").Value = 2440.718602890986
$ws.Range("R3").Value = 21966.46742601888
$ws.Range("S3").Value = 0.1423812039885396
$ws.Range("T3").Value = 0.1423812039885397

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 82.83048866666665
$ws.Range("H4").Value = 248.491466
$ws.Range("I4").Value = 0.3167437020391103
$ws.Range("J4").Value = 0.3167437020391103
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.31452
$ws.Range("N4").Value = 36.94356
$ws.Range("O4").Value = 0.1878601536983258
$ws.Range("P4").Value = 0.1878601536983257
$ws.Range("Q4").Value = 1020.01770929544
$ws.Range("R4").Value = 9180.159383658958
$ws.Range("S4").Value = 0.05950352054804395
$ws.Range("T4").Value = 0.05950352054804396

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 152.851481
$ws.Range("H5").Value = 458.554443
$ws.Range("I5").Value = 0.5845039034954311
$ws.Range("J5").Value = 0.5845039034954312
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 23.77057966666666
$ws.Range("N5").Value = 71.31173899999999
$ws.Range("O5").Value = 0.3626243450559418
$ws.Range("P5").Value = 0.3626243450559418
$ws.Range("Q5").Value = 3633.368306278486
$ws.Range("R5").Value = 32700.31475650637
$ws.Range("S5").Value = 0.2119553451876721
$ws.Range("T5").Value = 0.2119553451876721

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 152.851481
$ws.Range("H6").Value = 458.554443
$ws.Range("I6").Value = 0.5845039034954311
$ws.Range("J6").Value = 0.5845039034954312
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 29.46642766666666
$ws.Range("N6").Value = 88.399283
$ws.Range("O6").Value = 0.4495155012457325
$ws.Range("P6").Value = 0.4495155012457325
$ws.Range("Q6").Value = 4503.987108629374
$ws.Range("R6").Value = 40535.88397766437
$ws.Range("S6").Value = 0.262743565159836
$ws.Range("T6").Value = 0.262743565159836

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 152.851481
$ws.Range("H7").Value = 458.554443
$ws.Range("I7").Value = 0.5845039034954311
$ws.Range("J7").Value = 0.5845039034954312
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 12.31452
$ws.Range("N7").Value = 36.94356
$ws.Range("O7").Value = 0.1878601536983258
$ws.Range("P7").Value = 0.1878601536983257
$ws.Range("Q7").Value = 1882.29261980412
$ws.Range("R7").Value = 16940.63357823708
$ws.Range("S7").Value = 0.1098049931479231
$ws.Range("T7").Value = 0.1098049931479231

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 25.824378
$ws.Range("H8").Value = 77.47313399999999
$ws.Range("I8").Value = 0.09875239446545848
$ws.Range("J8").Value = 0.0987523944654585
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 23.77057966666666
$ws.Range("N8").Value = 71.31173899999999
$ws.Range("O8").Value = 0.3626243450559418
$ws.Range("P8").Value = 0.3626243450559418
$ws.Range("Q8").Value = 613.8604345911139
$ws.Range("R8").Value = 5524.743911320024
$ws.Range("S8").Value = 0.03581002236574289
$ws.Range("T8").Value = 0.0358100223657429

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 25.824378
$ws.Range("H9").Value = 77.47313399999999
$ws.Range("I9").Value = 0.09875239446545848
$ws.Range("J9").Value = 0.0987523944654585
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 29.46642766666666
$ws.Range("N9").Value = 88.399283
$ws.Range("O9").Value = 0.4495155012457325
$ws.Range("P9").Value = 0.4495155012457325
$ws.Range("Q9").Value = 760.9521663736579
$ws.Range("R9").Value = 6848.569497362921
$ws.Range("S9").Value = 0.04439073209735687
$ws.Range("T9").Value = 0.04439073209735688

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 25.824378
$ws.Range("H10").Value = 77.47313399999999
$ws.Range("I10").Value = 0.09875239446545848
$ws.Range("J10").Value = 0.0987523944654585
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 12.31452
$ws.Range("N10").Value = 36.94356
$ws.Range("O10").Value = 0.1878601536983258
$ws.Range("P10").Value = 0.1878601536983257
$ws.Range("Q10").Value = 318.0148193685599
$ws.Range("R10").Value = 2862.133374317039
$ws.Range("S10").Value = 0.01855164000235873
$ws.Range("T10").Value = 0.01855164000235873

